# The presentation has a single table on slide 5 (the table comparing
# "Type of document" / "Definition" / "Why it is important"). The author
# picked a different built-in table style for it from the Table Design
# gallery in PowerPoint, which rewrites the table's <a:tableStyleId>.
#
# Old style GUID: {60E381C7-C0AA-4806-9F8F-FB4DEF663BC5}
# New style GUID: {6C8CA814-BC42-46E6-BBA0-A0AA85ECB302}

$p = $ppt.ActivePresentation

# Slide 5 holds the table (shape 2 on that slide - a p:graphicFrame).
$slide = $p.Slides.Item(5)
$tableShape = $slide.Shapes.Item(2)

# Apply the new built-in table style by its style GUID, same as choosing
# a new style from the Table Design ribbon tab in PowerPoint.
$table = $tableShape.Table
$table.ApplyStyle("{6C8CA814-BC42-46E6-BBA0-A0AA85ECB302}")
